$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GH")

# Row 6 "Change in inventories" - update B6:F6
$ws.Range("B6").Value = -7535000.0
$ws.Range("C6").Value = -12971000.0
$ws.Range("D6").Value = -6178000.0
$ws.Range("E6").Value = -14892000.0
$ws.Range("F6").Value = -6045000.0

# Row 8 "Change in payables and accrued liability" - update B8:C8
$ws.Range("B8").Value = 65000000.0
$ws.Range("C8").Value = 74000000.0
